# NSMB.xlsx - "V4" sheet (ActiveSheet) edits
# Adds a 5-2 run's timing data, renames a couple of checkpoint labels, and
# moves the selection/scroll position further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22/24: update I/J timing data (K22/I24/J24 are formulas and will
#     recalculate automatically) ---
$ws.Range("I22").Value = 16243
$ws.Range("J22").Value = 18855

# --- Row 38 ("Get flag" / 5-1 area): new B38 entry start time so the
#     D38/F38 difference formulas produce a non-zero result ---
$ws.Range("B38").Value = 14776

# --- Row 42: relabel from "Enter 5-1" to the new "Enter 5-2" string, and add
#     the B42 start time ---
$ws.Range("A42").Value = "Enter 5-2"
$ws.Range("B42").Value = 15781

# --- Row 43: add B43 start time ---
$ws.Range("B43").Value = 16148

# --- Row 44: relabel from "Checkpoint 10" to the new "1st move" string, and
#     clear the G44/H44 checkpoint values (no longer applicable) ---
$ws.Range("A44").Value = "1st move"
$ws.Range("G44").ClearContents()
$ws.Range("H44").ClearContents()

# --- Row 54: add B54 start time ---
$ws.Range("B54").Value = 17646

# --- Row 55: add B55 start time ---
$ws.Range("B55").Value = 17913

# --- Move the current selection further down the sheet (author scrolled to
#     where they left off) ---
$ws.Range("B56").Select()
